# Auto-generated Excel COM-interop edit script.
# Applies the numeric corrections described by the commit's XML diff
# (per-leve currentAveragePrice / LevePrice / LeveProfit recompute).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 15
$ws.Range("H15").Value = 4551.5293
$ws.Range("I15").Value = 4551.5293
$ws.Range("K15").Value = 13654.5879
$ws.Range("M15").Value = -13485.5879

# row 62
$ws.Range("H62").Value = 12944.818
$ws.Range("I62").Value = 10299.125
$ws.Range("K62").Value = 10299.125
$ws.Range("M62").Value = -9675.125

# row 65
$ws.Range("H65").Value = 12944.818
$ws.Range("I65").Value = 10299.125
$ws.Range("K65").Value = 51495.625
$ws.Range("M65").Value = -48375.625

# row 87
$ws.Range("H87").Value = 99995
$ws.Range("J87").Value = 99995
$ws.Range("L87").Value = 99995
$ws.Range("N87").Value = -102491

# row 90
$ws.Range("H90").Value = 99995
$ws.Range("J90").Value = 99995
$ws.Range("L90").Value = 299985
$ws.Range("N90").Value = -312465

# row 125
$ws.Range("H125").Value = 3065.0715
$ws.Range("I125").Value = 1204
$ws.Range("K125").Value = 10836
$ws.Range("M125").Value = -8376

# row 134
$ws.Range("H134").Value = 90000
$ws.Range("J134").Value = 90000
$ws.Range("L134").Value = 90000
$ws.Range("N134").Value = -100140

# row 137
$ws.Range("H137").Value = 10824114
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 10824114
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 32472342
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -32477442

$ws = $wb.Worksheets.Item("ARM")
# row 2
$ws.Range("H2").Value = 8810.954
$ws.Range("I2").Value = 1795.2307
$ws.Range("J2").Value = 18944.777
$ws.Range("K2").Value = 1795.2307
$ws.Range("L2").Value = 18944.777
$ws.Range("M2").Value = -1682.2307
$ws.Range("N2").Value = -19170.777

# row 21
$ws.Range("H21").Value = 2321.25
$ws.Range("I21").Value = 262.5
$ws.Range("J21").Value = 8497.5
$ws.Range("K21").Value = 262.5
$ws.Range("L21").Value = 8497.5
$ws.Range("M21").Value = 111.5
$ws.Range("N21").Value = -9245.5

# row 30
$ws.Range("H30").Value = 1573.1111
$ws.Range("I30").Value = 894.75
$ws.Range("J30").Value = 7000
$ws.Range("K30").Value = 894.75
$ws.Range("L30").Value = 7000
$ws.Range("M30").Value = -744.75
$ws.Range("N30").Value = -7300

# row 32
$ws.Range("H32").Value = 5467117
$ws.Range("I32").Value = 5749894
$ws.Range("K32").Value = 5749894
$ws.Range("M32").Value = -5749607

# row 61
$ws.Range("H61").Value = 1392758.1
$ws.Range("I61").Value = 1855051.5
$ws.Range("K61").Value = 1855051.5
$ws.Range("M61").Value = -1854839.5

# row 97
$ws.Range("H97").Value = 1070.6364
$ws.Range("I97").Value = 879.55554
$ws.Range("J97").Value = 1930.5
$ws.Range("K97").Value = 879.55554
$ws.Range("L97").Value = 1930.5
$ws.Range("M97").Value = -383.55554
$ws.Range("N97").Value = -2922.5

# row 116
$ws.Range("H116").Value = 8810.954
$ws.Range("I116").Value = 1795.2307
$ws.Range("J116").Value = 18944.777
$ws.Range("K116").Value = 1795.2307
$ws.Range("L116").Value = 18944.777
$ws.Range("M116").Value = 498.7692999999999
$ws.Range("N116").Value = -23532.777

# row 136
$ws.Range("H136").Value = 1392758.1
$ws.Range("I136").Value = 1855051.5
$ws.Range("K136").Value = 5565154.5
$ws.Range("M136").Value = -5562604.5

$ws = $wb.Worksheets.Item("BSM")
# row 3
$ws.Range("H3").Value = 8810.954
$ws.Range("I3").Value = 1795.2307
$ws.Range("J3").Value = 18944.777
$ws.Range("K3").Value = 1795.2307
$ws.Range("L3").Value = 18944.777
$ws.Range("M3").Value = -1681.2307
$ws.Range("N3").Value = -19172.777

# row 107
$ws.Range("H107").Value = 8830.799999999999
$ws.Range("I107").Value = 9367.666999999999
$ws.Range("K107").Value = 9367.666999999999
$ws.Range("M107").Value = -7447.666999999999

# row 124
$ws.Range("H124").Value = 25000
$ws.Range("I124").Value = 25000
$ws.Range("K124").Value = 25000
$ws.Range("M124").Value = -20090

$ws = $wb.Worksheets.Item("CRP")
# row 4
$ws.Range("H4").Value = 40000000
$ws.Range("I4").Value = 40000000
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 40000000
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -39999888
$ws.Range("N4").ClearContents()

# row 31
$ws.Range("H31").Value = 108274.06
$ws.Range("J31").Value = 23892.334
$ws.Range("L31").Value = 23892.334
$ws.Range("N31").Value = -24482.334

# row 34
$ws.Range("H34").Value = 108274.06
$ws.Range("J34").Value = 23892.334
$ws.Range("L34").Value = 23892.334
$ws.Range("N34").Value = -24296.334

# row 41
$ws.Range("H41").Value = 9008.429
$ws.Range("J41").Value = 10000
$ws.Range("L41").Value = 10000
$ws.Range("N41").Value = -10856

# row 94
$ws.Range("H94").Value = 2153.3
$ws.Range("J94").Value = 2717.4285
$ws.Range("L94").Value = 2717.4285
$ws.Range("N94").Value = -3619.4285

# row 132
$ws.Range("H132").Value = 34663164
$ws.Range("I132").Value = 58826396
$ws.Range("J132").Value = 431918.1
$ws.Range("K132").Value = 176479188
$ws.Range("L132").Value = 1295754.3
$ws.Range("M132").Value = -176476658
$ws.Range("N132").Value = -1300814.3

$ws = $wb.Worksheets.Item("CUL")
# row 5
$ws.Range("H5").Value = 596.4286
$ws.Range("I5").Value = 444
$ws.Range("J5").Value = 977.5
$ws.Range("K5").Value = 1332
$ws.Range("L5").Value = 2932.5
$ws.Range("M5").Value = -1220
$ws.Range("N5").Value = -3156.5

# row 92
$ws.Range("H92").Value = 501.83334
$ws.Range("I92").Value = 501.83334
$ws.Range("K92").Value = 1505.50002
$ws.Range("M92").Value = -257.5000199999999

# row 135
$ws.Range("H135").Value = 596.4286
$ws.Range("I135").Value = 444
$ws.Range("J135").Value = 977.5
$ws.Range("K135").Value = 3996
$ws.Range("L135").Value = 8797.5
$ws.Range("M135").Value = -1461
$ws.Range("N135").Value = -13867.5

# row 137
$ws.Range("H137").Value = 3778
$ws.Range("I137").Value = 2291.5715
$ws.Range("J137").Value = 7246.3335
$ws.Range("K137").Value = 6874.7145
$ws.Range("L137").Value = 21739.0005
$ws.Range("M137").Value = -1774.7145
$ws.Range("N137").Value = -31939.0005

# row 140
$ws.Range("H140").Value = 2747
$ws.Range("I140").Value = 1869.9
$ws.Range("K140").Value = 5609.700000000001
$ws.Range("M140").Value = -429.7000000000007

$ws = $wb.Worksheets.Item("GSM")
# row 102
$ws.Range("H102").Value = 3683.7144
$ws.Range("J102").Value = 7499.6
$ws.Range("L102").Value = 7499.6
$ws.Range("N102").Value = -10743.6

# row 107
$ws.Range("H107").Value = 33432.727
$ws.Range("J107").Value = 2220.6
$ws.Range("L107").Value = 2220.6
$ws.Range("N107").Value = -6060.6

# row 122
$ws.Range("H122").Value = 56618.95
$ws.Range("I122").Value = 114592.664
$ws.Range("K122").Value = 343777.992
$ws.Range("M122").Value = -341327.992

# row 132
$ws.Range("H132").Value = 27357060
$ws.Range("I132").Value = 42171756
$ws.Range("J132").Value = 6846.615
$ws.Range("K132").Value = 126515268
$ws.Range("L132").Value = 20539.845
$ws.Range("M132").Value = -126512738
$ws.Range("N132").Value = -25599.845

$ws = $wb.Worksheets.Item("LTW")
# row 7
$ws.Range("H7").Value = 3999.9443
$ws.Range("I7").Value = 3884.5386
$ws.Range("K7").Value = 3884.5386
$ws.Range("M7").Value = -3772.5386

# row 115
$ws.Range("H115").Value = 78000
$ws.Range("J115").Value = 78000
$ws.Range("L115").Value = 78000
$ws.Range("N115").Value = -80350

# row 126
$ws.Range("H126").Value = 3999.9443
$ws.Range("I126").Value = 3884.5386
$ws.Range("K126").Value = 11653.6158
$ws.Range("M126").Value = -9183.6158

# row 132
$ws.Range("H132").Value = 4978230.5
$ws.Range("I132").Value = 11602563
$ws.Range("J132").Value = 9981.25
$ws.Range("K132").Value = 34807689
$ws.Range("L132").Value = 29943.75
$ws.Range("M132").Value = -34805159
$ws.Range("N132").Value = -35003.75

$ws = $wb.Worksheets.Item("WVR")
# row 95
$ws.Range("H95").Value = 51580.57
$ws.Range("J95").Value = 51580.57
$ws.Range("L95").Value = 51580.57
$ws.Range("N95").Value = -57072.57

# row 131
$ws.Range("H131").Value = 80000
$ws.Range("J131").Value = 80000
$ws.Range("L131").Value = 80000
$ws.Range("N131").Value = -90080

# row 132
$ws.Range("H132").Value = 5597307
$ws.Range("I132").Value = 8753287
$ws.Range("J132").Value = 13649.077
$ws.Range("K132").Value = 26259861
$ws.Range("L132").Value = 40947.231
$ws.Range("M132").Value = -26257331

# row 136
$ws.Range("H136").Value = 20466.059
$ws.Range("I136").Value = 13986.667
$ws.Range("K136").Value = 41960.001
$ws.Range("M136").Value = -39410.001

